$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # 展览
$ws.Cells.Item(3, 6).Value = 985
$ws.Cells.Item(4, 6).Value = 618
$ws.Cells.Item(5, 6).Value = 3097
$ws.Cells.Item(6, 6).Value = 828
$ws.Cells.Item(7, 6).Value = 607
$ws.Cells.Item(8, 6).Value = 606
$ws.Cells.Item(9, 6).Value = 473
$ws.Cells.Item(12, 6).Value = 612
$ws.Cells.Item(14, 6).Value = 2172
$ws.Cells.Item(15, 6).Value = 1270
$ws.Cells.Item(16, 6).Value = 756
$ws.Cells.Item(19, 6).Value = 2698
$ws.Cells.Item(23, 6).Value = 537
$ws.Cells.Item(24, 6).Value = 517
$ws.Cells.Item(25, 6).Value = 726
$ws.Cells.Item(26, 6).Value = 35
$ws.Cells.Item(27, 6).Value = 30
$ws.Cells.Item(29, 6).Value = 32
$ws.Cells.Item(30, 6).Value = 566
$ws.Cells.Item(31, 6).Value = 570
$ws.Cells.Item(33, 6).Value = 132
$ws.Cells.Item(34, 6).Value = 920
$ws.Cells.Item(35, 6).Value = 4743
$ws.Cells.Item(36, 6).Value = 294
$ws.Cells.Item(37, 6).Value = 59
$ws.Cells.Item(38, 6).Value = 24
$ws.Cells.Item(39, 6).Value = 85
$ws = $wb.Worksheets.Item(2)  # 演出
$ws.Cells.Item(8, 6).Value = 373
$ws.Cells.Item(22, 6).Value = 9
$ws.Cells.Item(23, 6).Value = 282
$ws.Cells.Item(24, 6).Value = 28
$ws.Cells.Item(25, 6).Value = 11
$ws.Cells.Item(26, 6).Value = 313
$ws.Cells.Item(28, 6).Value = 351
$ws.Cells.Item(32, 6).Value = 40
$ws.Cells.Item(38, 6).Value = 628
$ws.Cells.Item(39, 6).Value = 628
$ws.Cells.Item(40, 6).Value = 23
$ws = $wb.Worksheets.Item(3)  # 本地生活
$ws.Cells.Item(4, 6).Value = 1480
$ws.Cells.Item(5, 6).Value = 585
$ws.Cells.Item(6, 6).Value = 309
$ws.Cells.Item(7, 6).Value = 290
$ws = $wb.Worksheets.Item(4)  # 全部类型
$ws.Cells.Item(3, 6).Value = 1480
$ws.Cells.Item(4, 6).Value = 585
$ws.Cells.Item(6, 6).Value = 309
$ws.Cells.Item(7, 6).Value = 985
$ws.Cells.Item(8, 6).Value = 618
$ws.Cells.Item(9, 6).Value = 3097
$ws.Cells.Item(10, 6).Value = 828
$ws.Cells.Item(11, 6).Value = 607
$ws.Cells.Item(12, 6).Value = 606
$ws.Cells.Item(13, 6).Value = 473
$ws.Cells.Item(16, 6).Value = 612
$ws.Cells.Item(21, 6).Value = 2172
$ws.Cells.Item(22, 6).Value = 1270
$ws.Cells.Item(23, 6).Value = 756
$ws.Cells.Item(24, 6).Value = 47
$ws.Cells.Item(26, 6).Value = 2698
$ws.Cells.Item(29, 6).Value = 537
$ws.Cells.Item(31, 6).Value = 290
$ws.Cells.Item(33, 6).Value = 517
$ws.Cells.Item(34, 6).Value = 726
$ws.Cells.Item(35, 6).Value = 726
$ws.Cells.Item(36, 6).Value = 35
$ws.Cells.Item(37, 6).Value = 9
$ws.Cells.Item(38, 6).Value = 32
$ws.Cells.Item(39, 6).Value = 28
$ws.Cells.Item(40, 6).Value = 11
$ws.Cells.Item(41, 6).Value = 566
$ws.Cells.Item(42, 6).Value = 313
$ws.Cells.Item(44, 6).Value = 920
$ws.Cells.Item(45, 6).Value = 4743
$ws.Cells.Item(46, 6).Value = 40
$ws.Cells.Item(47, 6).Value = 294
$ws.Cells.Item(49, 6).Value = 59
$ws.Cells.Item(51, 6).Value = 628
